$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new Gender column, matching style of existing headers (H1/I1/J1)
$ws.Range("K1").Value = "Gender"
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)

# Fill formulas down K2:K110. Mirror the existing J-column shared-formula
# grouping (standalone J2, shared J3:J66, shared J67:J110) by writing the
# new column's formulas in the same three blocks.
$ws.Range("K2").Formula = '=IF(LEFT(D2,1)="M","M","F")'
$ws.Range("K3:K66").Formula = '=IF(LEFT(D3,1)="M","M","F")'
$ws.Range("K67:K110").Formula = '=IF(LEFT(D67,1)="M","M","F")'

# Update selection to K2, clear the scrolled top-left cell
$ws.Application.Goto($ws.Range("K2"))
$ws.Range("K2").Select()
